$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 279; this shifts the existing rows 279-287
# down to 280-288, matching the diff (old row N becomes row N+1).
$ws.Rows.Item(279).Insert()

# Populate the newly inserted row 279 with the new record.
$ws.Cells.Item(279, 1).Value = 10
$ws.Cells.Item(279, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(279, 3).Value = "La Araucanía"
$ws.Cells.Item(279, 4).Value = 45075
$ws.Cells.Item(279, 4).Style = $ws.Cells.Item(280, 4).Style
$ws.Cells.Item(279, 4).NumberFormat = $ws.Cells.Item(280, 4).NumberFormat
$ws.Cells.Item(279, 5).Value = 9
$ws.Cells.Item(279, 6).Value = 100112013
$ws.Cells.Item(279, 7).Value = "Alcachofa"
$ws.Cells.Item(279, 8).Value = "Madrigal"
$ws.Cells.Item(279, 9).Value = "Primera"
$ws.Cells.Item(279, 10).Value = 55
$ws.Cells.Item(279, 11).Value = 15000
$ws.Cells.Item(279, 12).Value = 15000
$ws.Cells.Item(279, 13).Value = 15000
$ws.Cells.Item(279, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(279, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(279, 16).Value = 375
$ws.Cells.Item(279, 17).Value = 40
$ws.Cells.Item(279, 18).Value = "Hortaliza"
